$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 5: 保險 (insurance)
# Fix the header row (B1:D1 previously held stray data values instead of the
# actual column headers) and append the standard metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that every other sheet already carries.
# ---------------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(5)

# Header row
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Make header cells bold / centered / bordered, matching the rest of row 1
$hdrRangeIns = $wsIns.Range("B1:K1")
$hdrRangeIns.Font.Bold = $true
$hdrRangeIns.HorizontalAlignment = -4108
$hdrRangeIns.VerticalAlignment = -4160
$hdrRangeIns.Borders.LineStyle = 1

# New metadata columns for each data row (values identical across rows)
$insRows = 75,76,77,78
for ($i = 0; $i -lt $insRows.Length; $i++) {
    $r = $i + 2
    $wsIns.Cells.Item($r, 5).Value = "insurance"
    $wsIns.Cells.Item($r, 6).Value = "normal"
    $wsIns.Cells.Item($r, 7).NumberFormat = "@"
    $wsIns.Cells.Item($r, 7).Value = "2012-04-12"
    $wsIns.Cells.Item($r, 8).Value = "王育敏"
    $wsIns.Cells.Item($r, 9).Value = 1728
    $wsIns.Cells.Item($r, 10).Value = "tmp48bc1"
    $wsIns.Cells.Item($r, 11).Value = $insRows[$i]
}

# ---------------------------------------------------------------------------
# Sheet 6: 債務 (debt)
# Same treatment: fix the header row and append the standard metadata
# columns.
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(6)

# Header row
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

$hdrRangeDebt = $wsDebt.Range("B1:N1")
$hdrRangeDebt.Font.Bold = $true
$hdrRangeDebt.HorizontalAlignment = -4108
$hdrRangeDebt.VerticalAlignment = -4160
$hdrRangeDebt.Borders.LineStyle = 1

# New metadata columns for row 2
$wsDebt.Cells.Item(2, 8).Value = "debt"
$wsDebt.Cells.Item(2, 9).Value = "normal"
$wsDebt.Cells.Item(2, 10).NumberFormat = "@"
$wsDebt.Cells.Item(2, 10).Value = "2012-04-12"
$wsDebt.Cells.Item(2, 11).Value = "王育敏"
$wsDebt.Cells.Item(2, 12).Value = 1728
$wsDebt.Cells.Item(2, 13).Value = "tmp48bc1"
$wsDebt.Cells.Item(2, 14).Value = 88
